$d = $word.ActiveDocument

$d.Content.Find.Execute("59+0=", $true, $false, $false, $false, $false, $true, 1, $false, "50-38=", 2) | Out-Null
$d.Content.Find.Execute("9+55=", $true, $false, $false, $false, $false, $true, 1, $false, "32+18=", 2) | Out-Null
$d.Content.Find.Execute("53+13=", $true, $false, $false, $false, $false, $true, 1, $false, "77-59=", 2) | Out-Null
$d.Content.Find.Execute("54-9=", $true, $false, $false, $false, $false, $true, 1, $false, "66+11=", 2) | Out-Null
$d.Content.Find.Execute("45-16=", $true, $false, $false, $false, $false, $true, 1, $false, "39+60=", 2) | Out-Null
$d.Content.Find.Execute("49+32=", $true, $false, $false, $false, $false, $true, 1, $false, "13+25=", 2) | Out-Null
$d.Content.Find.Execute("17-10=", $true, $false, $false, $false, $false, $true, 1, $false, "53+1=", 2) | Out-Null
$d.Content.Find.Execute("96-66=", $true, $false, $false, $false, $false, $true, 1, $false, "27+16=", 2) | Out-Null
$d.Content.Find.Execute("53-40=", $true, $false, $false, $false, $false, $true, 1, $false, "39-17=", 2) | Out-Null
$d.Content.Find.Execute("2+4=", $true, $false, $false, $false, $false, $true, 1, $false, "87-35=", 2) | Out-Null
$d.Content.Find.Execute("40+36=", $true, $false, $false, $false, $false, $true, 1, $false, "18-0=", 2) | Out-Null
$d.Content.Find.Execute("8+45=", $true, $false, $false, $false, $false, $true, 1, $false, "24-11=", 2) | Out-Null
$d.Content.Find.Execute("72-48=", $true, $false, $false, $false, $false, $true, 1, $false, "13+8=", 2) | Out-Null
$d.Content.Find.Execute("94-85=", $true, $false, $false, $false, $false, $true, 1, $false, "19+8=", 2) | Out-Null
$d.Content.Find.Execute("90-1=", $true, $false, $false, $false, $false, $true, 1, $false, "82-69=", 2) | Out-Null
$d.Content.Find.Execute("41+39=", $true, $false, $false, $false, $false, $true, 1, $false, "96-40=", 2) | Out-Null
$d.Content.Find.Execute("34+33=", $true, $false, $false, $false, $false, $true, 1, $false, "79+19=", 2) | Out-Null
$d.Content.Find.Execute("2+72=", $true, $false, $false, $false, $false, $true, 1, $false, "88-30=", 2) | Out-Null
$d.Content.Find.Execute("7+24=", $true, $false, $false, $false, $false, $true, 1, $false, "70+27=", 2) | Out-Null
$d.Content.Find.Execute("79-15=", $true, $false, $false, $false, $false, $true, 1, $false, "96-42=", 2) | Out-Null
$d.Content.Find.Execute("67-2=", $true, $false, $false, $false, $false, $true, 1, $false, "98-1=", 2) | Out-Null
$d.Content.Find.Execute("95-22=", $true, $false, $false, $false, $false, $true, 1, $false, "73-10=", 2) | Out-Null
$d.Content.Find.Execute("40-18=", $true, $false, $false, $false, $false, $true, 1, $false, "15-15=", 2) | Out-Null
$d.Content.Find.Execute("0+91=", $true, $false, $false, $false, $false, $true, 1, $false, "5+32=", 2) | Out-Null
$d.Content.Find.Execute("73-37=", $true, $false, $false, $false, $false, $true, 1, $false, "16+15=", 2) | Out-Null
$d.Content.Find.Execute("68+15=", $true, $false, $false, $false, $false, $true, 1, $false, "57-10=", 2) | Out-Null
$d.Content.Find.Execute("36+8=", $true, $false, $false, $false, $false, $true, 1, $false, "56-1=", 2) | Out-Null
$d.Content.Find.Execute("26+72=", $true, $false, $false, $false, $false, $true, 1, $false, "43+55=", 2) | Out-Null
$d.Content.Find.Execute("21-8=", $true, $false, $false, $false, $false, $true, 1, $false, "38-19=", 2) | Out-Null
$d.Content.Find.Execute("52+26=", $true, $false, $false, $false, $false, $true, 1, $false, "25+72=", 2) | Out-Null
$d.Content.Find.Execute("69-17=", $true, $false, $false, $false, $false, $true, 1, $false, "34-12=", 2) | Out-Null
$d.Content.Find.Execute("23+60=", $true, $false, $false, $false, $false, $true, 1, $false, "18-9=", 2) | Out-Null
$d.Content.Find.Execute("29+39=", $true, $false, $false, $false, $false, $true, 1, $false, "47-45=", 2) | Out-Null
$d.Content.Find.Execute("18+58=", $true, $false, $false, $false, $false, $true, 1, $false, "55-10=", 2) | Out-Null
$d.Content.Find.Execute("98-23=", $true, $false, $false, $false, $false, $true, 1, $false, "95-82=", 2) | Out-Null
$d.Content.Find.Execute("20-16=", $true, $false, $false, $false, $false, $true, 1, $false, "69-51=", 2) | Out-Null
$d.Content.Find.Execute("78-45=", $true, $false, $false, $false, $false, $true, 1, $false, "42-19=", 2) | Out-Null
$d.Content.Find.Execute("99-54=", $true, $false, $false, $false, $false, $true, 1, $false, "76-29=", 2) | Out-Null
$d.Content.Find.Execute("46+22=", $true, $false, $false, $false, $false, $true, 1, $false, "81-76=", 2) | Out-Null
$d.Content.Find.Execute("75-37=", $true, $false, $false, $false, $false, $true, 1, $false, "61+29=", 2) | Out-Null
$d.Content.Find.Execute("20+58=", $true, $false, $false, $false, $false, $true, 1, $false, "4+80=", 2) | Out-Null
$d.Content.Find.Execute("73-24=", $true, $false, $false, $false, $false, $true, 1, $false, "66+29=", 2) | Out-Null
$d.Content.Find.Execute("47-39=", $true, $false, $false, $false, $false, $true, 1, $false, "97-63=", 2) | Out-Null
$d.Content.Find.Execute("89-73=", $true, $false, $false, $false, $false, $true, 1, $false, "68+30=", 2) | Out-Null
$d.Content.Find.Execute("84-65=", $true, $false, $false, $false, $false, $true, 1, $false, "5+71=", 2) | Out-Null
$d.Content.Find.Execute("2+12=", $true, $false, $false, $false, $false, $true, 1, $false, "26+71=", 2) | Out-Null
$d.Content.Find.Execute("73-44=", $true, $false, $false, $false, $false, $true, 1, $false, "72+25=", 2) | Out-Null
$d.Content.Find.Execute("91-15=", $true, $false, $false, $false, $false, $true, 1, $false, "73+11=", 2) | Out-Null
$d.Content.Find.Execute("56-43=", $true, $false, $false, $false, $false, $true, 1, $false, "85-28=", 2) | Out-Null
$d.Content.Find.Execute("87-51=", $true, $false, $false, $false, $false, $true, 1, $false, "48-18=", 2) | Out-Null
$d.Content.Find.Execute("8-8=", $true, $false, $false, $false, $false, $true, 1, $false, "81+5=", 2) | Out-Null
$d.Content.Find.Execute("20+34=", $true, $false, $false, $false, $false, $true, 1, $false, "61+0=", 2) | Out-Null
$d.Content.Find.Execute("92-76=", $true, $false, $false, $false, $false, $true, 1, $false, "44-10=", 2) | Out-Null
$d.Content.Find.Execute("14-0=", $true, $false, $false, $false, $false, $true, 1, $false, "54+31=", 2) | Out-Null
$d.Content.Find.Execute("55+14=", $true, $false, $false, $false, $false, $true, 1, $false, "94+0=", 2) | Out-Null
$d.Content.Find.Execute("77-11=", $true, $false, $false, $false, $false, $true, 1, $false, "10+28=", 2) | Out-Null
$d.Content.Find.Execute("46+30=", $true, $false, $false, $false, $false, $true, 1, $false, "57-14=", 2) | Out-Null
$d.Content.Find.Execute("67-58=", $true, $false, $false, $false, $false, $true, 1, $false, "58+26=", 2) | Out-Null
$d.Content.Find.Execute("62+31=", $true, $false, $false, $false, $false, $true, 1, $false, "97-70=", 2) | Out-Null
$d.Content.Find.Execute("46-32=", $true, $false, $false, $false, $false, $true, 1, $false, "47-41=", 2) | Out-Null
$d.Content.Find.Execute("84-17=", $true, $false, $false, $false, $false, $true, 1, $false, "2+58=", 2) | Out-Null
$d.Content.Find.Execute("29-4=", $true, $false, $false, $false, $false, $true, 1, $false, "49-47=", 2) | Out-Null
$d.Content.Find.Execute("53+30=", $true, $false, $false, $false, $false, $true, 1, $false, "16-15=", 2) | Out-Null
$d.Content.Find.Execute("83-27=", $true, $false, $false, $false, $false, $true, 1, $false, "26+61=", 2) | Out-Null
$d.Content.Find.Execute("40-31=", $true, $false, $false, $false, $false, $true, 1, $false, "84-59=", 2) | Out-Null
$d.Content.Find.Execute("0+17=", $true, $false, $false, $false, $false, $true, 1, $false, "57-46=", 2) | Out-Null
$d.Content.Find.Execute("83-35=", $true, $false, $false, $false, $false, $true, 1, $false, "90-6=", 2) | Out-Null
$d.Content.Find.Execute("0+95=", $true, $false, $false, $false, $false, $true, 1, $false, "79+20=", 2) | Out-Null
$d.Content.Find.Execute("71-39=", $true, $false, $false, $false, $false, $true, 1, $false, "56+24=", 2) | Out-Null
$d.Content.Find.Execute("93-37=", $true, $false, $false, $false, $false, $true, 1, $false, "58-50=", 2) | Out-Null
$d.Content.Find.Execute("18+15=", $true, $false, $false, $false, $false, $true, 1, $false, "27-7=", 2) | Out-Null
$d.Content.Find.Execute("91-40=", $true, $false, $false, $false, $false, $true, 1, $false, "46+51=", 2) | Out-Null
$d.Content.Find.Execute("81+1=", $true, $false, $false, $false, $false, $true, 1, $false, "10-10=", 2) | Out-Null
$d.Content.Find.Execute("81-43=", $true, $false, $false, $false, $false, $true, 1, $false, "77-6=", 2) | Out-Null
$d.Content.Find.Execute("8+90=", $true, $false, $false, $false, $false, $true, 1, $false, "33+51=", 2) | Out-Null
$d.Content.Find.Execute("88+9=", $true, $false, $false, $false, $false, $true, 1, $false, "81-19=", 2) | Out-Null
$d.Content.Find.Execute("31+63=", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=", 2) | Out-Null
$d.Content.Find.Execute("68-67=", $true, $false, $false, $false, $false, $true, 1, $false, "26+41=", 2) | Out-Null
$d.Content.Find.Execute("64-44=", $true, $false, $false, $false, $false, $true, 1, $false, "46-9=", 2) | Out-Null
$d.Content.Find.Execute("46-19=", $true, $false, $false, $false, $false, $true, 1, $false, "56+40=", 2) | Out-Null
$d.Content.Find.Execute("9+15=", $true, $false, $false, $false, $false, $true, 1, $false, "39+36=", 2) | Out-Null
$d.Content.Find.Execute("32-14=", $true, $false, $false, $false, $false, $true, 1, $false, "68+9=", 2) | Out-Null
$d.Content.Find.Execute("42-20=", $true, $false, $false, $false, $false, $true, 1, $false, "14-7=", 2) | Out-Null
$d.Content.Find.Execute("31+6=", $true, $false, $false, $false, $false, $true, 1, $false, "93-5=", 2) | Out-Null
$d.Content.Find.Execute("61-61=", $true, $false, $false, $false, $false, $true, 1, $false, "18+48=", 2) | Out-Null
$d.Content.Find.Execute("60-16=", $true, $false, $false, $false, $false, $true, 1, $false, "11-2=", 2) | Out-Null
$d.Content.Find.Execute("92-1=", $true, $false, $false, $false, $false, $true, 1, $false, "96-29=", 2) | Out-Null
$d.Content.Find.Execute("85-22=", $true, $false, $false, $false, $false, $true, 1, $false, "69+15=", 2) | Out-Null
$d.Content.Find.Execute("95-69=", $true, $false, $false, $false, $false, $true, 1, $false, "72-70=", 2) | Out-Null
$d.Content.Find.Execute("95-48=", $true, $false, $false, $false, $false, $true, 1, $false, "55-16=", 2) | Out-Null
$d.Content.Find.Execute("17-16=", $true, $false, $false, $false, $false, $true, 1, $false, "83-36=", 2) | Out-Null
$d.Content.Find.Execute("35+56=", $true, $false, $false, $false, $false, $true, 1, $false, "60-37=", 2) | Out-Null
$d.Content.Find.Execute("86-72=", $true, $false, $false, $false, $false, $true, 1, $false, "7+90=", 2) | Out-Null
$d.Content.Find.Execute("52+28=", $true, $false, $false, $false, $false, $true, 1, $false, "96-54=", 2) | Out-Null
$d.Content.Find.Execute("86-40=", $true, $false, $false, $false, $false, $true, 1, $false, "50+6=", 2) | Out-Null
$d.Content.Find.Execute("79-59=", $true, $false, $false, $false, $false, $true, 1, $false, "78-71=", 2) | Out-Null
$d.Content.Find.Execute("9+3=", $true, $false, $false, $false, $false, $true, 1, $false, "35+60=", 2) | Out-Null
$d.Content.Find.Execute("91-62=", $true, $false, $false, $false, $false, $true, 1, $false, "79-4=", 2) | Out-Null
$d.Content.Find.Execute("42-15=", $true, $false, $false, $false, $false, $true, 1, $false, "87-58=", 2) | Out-Null
$d.Content.Find.Execute("82-20=", $true, $false, $false, $false, $false, $true, 1, $false, "29-27=", 2) | Out-Null
